# Feria Lagunitas de Puerto Montt - Choclo: add 4 new weekly records.
# This inserts 4 rows at row 176 (pushing the existing rows 176-187 down
# to 180-191, unchanged) and fills the newly inserted rows 176-179 with
# the new "Choclero"/"Dulce o Americano" price records for O'Higgins /
# Región Metropolitana reported on 44585 (2022-01-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 176 (shifts 176:187 -> 180:191)
$ws.Range("A176:A179").EntireRow.Insert()

# New row 176
$ws.Cells.Item(176,1).Value  = 4
$ws.Cells.Item(176,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(176,3).Value  = "Los Lagos"
$ws.Cells.Item(176,4).Value  = 44585
$ws.Cells.Item(176,5).Value  = 10
$ws.Cells.Item(176,6).Value  = 100112024
$ws.Cells.Item(176,7).Value  = "Choclo"
$ws.Cells.Item(176,8).Value  = "Choclero"
$ws.Cells.Item(176,9).Value  = "Primera"
$ws.Cells.Item(176,10).Value = 5000
$ws.Cells.Item(176,11).Value = 500
$ws.Cells.Item(176,12).Value = 500
$ws.Cells.Item(176,13).Value = 500
$ws.Cells.Item(176,14).Value = "$/unidad"
$ws.Cells.Item(176,15).Value = "Región de O'Higgins"
$ws.Cells.Item(176,16).Value = 500
$ws.Cells.Item(176,17).Value = 1
$ws.Cells.Item(176,18).Value = "Hortaliza"

# New row 177
$ws.Cells.Item(177,1).Value  = 4
$ws.Cells.Item(177,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(177,3).Value  = "Los Lagos"
$ws.Cells.Item(177,4).Value  = 44585
$ws.Cells.Item(177,5).Value  = 10
$ws.Cells.Item(177,6).Value  = 100112024
$ws.Cells.Item(177,7).Value  = "Choclo"
$ws.Cells.Item(177,8).Value  = "Choclero"
$ws.Cells.Item(177,9).Value  = "Segunda"
$ws.Cells.Item(177,10).Value = 5000
$ws.Cells.Item(177,11).Value = 300
$ws.Cells.Item(177,12).Value = 300
$ws.Cells.Item(177,13).Value = 300
$ws.Cells.Item(177,14).Value = "$/unidad"
$ws.Cells.Item(177,15).Value = "Región de O'Higgins"
$ws.Cells.Item(177,16).Value = 300
$ws.Cells.Item(177,17).Value = 1
$ws.Cells.Item(177,18).Value = "Hortaliza"

# New row 178
$ws.Cells.Item(178,1).Value  = 4
$ws.Cells.Item(178,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(178,3).Value  = "Los Lagos"
$ws.Cells.Item(178,4).Value  = 44585
$ws.Cells.Item(178,5).Value  = 10
$ws.Cells.Item(178,6).Value  = 100112024
$ws.Cells.Item(178,7).Value  = "Choclo"
$ws.Cells.Item(178,8).Value  = "Dulce o Americano"
$ws.Cells.Item(178,9).Value  = "Primera"
$ws.Cells.Item(178,10).Value = 10000
$ws.Cells.Item(178,11).Value = 250
$ws.Cells.Item(178,12).Value = 250
$ws.Cells.Item(178,13).Value = 250
$ws.Cells.Item(178,14).Value = "$/unidad"
$ws.Cells.Item(178,15).Value = "Región Metropolitana"
$ws.Cells.Item(178,16).Value = 250
$ws.Cells.Item(178,17).Value = 1
$ws.Cells.Item(178,18).Value = "Hortaliza"

# New row 179
$ws.Cells.Item(179,1).Value  = 4
$ws.Cells.Item(179,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(179,3).Value  = "Los Lagos"
$ws.Cells.Item(179,4).Value  = 44585
$ws.Cells.Item(179,5).Value  = 10
$ws.Cells.Item(179,6).Value  = 100112024
$ws.Cells.Item(179,7).Value  = "Choclo"
$ws.Cells.Item(179,8).Value  = "Dulce o Americano"
$ws.Cells.Item(179,9).Value  = "Segunda"
$ws.Cells.Item(179,10).Value = 10000
$ws.Cells.Item(179,11).Value = 200
$ws.Cells.Item(179,12).Value = 200
$ws.Cells.Item(179,13).Value = 200
$ws.Cells.Item(179,14).Value = "$/unidad"
$ws.Cells.Item(179,15).Value = "Región Metropolitana"
$ws.Cells.Item(179,16).Value = 200
$ws.Cells.Item(179,17).Value = 1
$ws.Cells.Item(179,18).Value = "Hortaliza"
